$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row labels to match new naming convention
# Order matches the order new shared strings were appended in the target file
$ws.Range("E1").Value = "description"
$ws.Range("F1").Value = "soils and landform"
$ws.Range("D1").Value = "notes"
$ws.Range("B1").Value = "longitude (deg)"
$ws.Range("C1").Value = "latitude (deg)"

# Update selection to B2 (cursor position as recorded in saved file)
$ws.Range("B2").Select()
